# Apply "Natmi following Dr Hou advice" edit:
# Rebuild the LR-pair result rows (rows 2-7) into the full cross-join of
# sending/target clusters [ECs, FAPs, sCs] x [ECs, FAPs, M2, sCs] (rows 2-13),
# with updated statistics for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data rows so the sheet can be rebuilt with the new 12-row table
$ws.Range("A2:T7").ClearContents()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rarres2"
$ws.Range("C2").Value = "Ccrl2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.9064200000000001
$ws.Range("H2").Value = 2.71926
$ws.Range("I2").Value = 0.01741933661286065
$ws.Range("J2").Value = 0.01741933661286065
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 169.8267096666667
$ws.Range("N2").Value = 509.480129
$ws.Range("O2").Value = 0.7951622771939805
$ws.Range("P2").Value = 0.7951622771939802
$ws.Range("Q2").Value = 153.93432617606
$ws.Range("R2").Value = 1385.40893558454
$ws.Range("S2").Value = 0.01385119936829075
$ws.Range("T2").Value = 0.01385119936829075

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rarres2"
$ws.Range("C3").Value = "Ccrl2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.9064200000000001
$ws.Range("H3").Value = 2.71926
$ws.Range("I3").Value = 0.01741933661286065
$ws.Range("J3").Value = 0.01741933661286065
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.30056
$ws.Range("N3").Value = 0.90168
$ws.Range("O3").Value = 0.001407281425298274
$ws.Range("P3").Value = 0.001407281425298273
$ws.Range("Q3").Value = 0.2724335952
$ws.Range("R3").Value = 2.4519023568
$ws.Range("S3").Value = [double]"2.451390885629694e-05"
$ws.Range("T3").Value = [double]"2.451390885629694e-05"

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Rarres2"
$ws.Range("C4").Value = "Ccrl2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.9064200000000001
$ws.Range("H4").Value = 2.71926
$ws.Range("I4").Value = 0.01741933661286065
$ws.Range("J4").Value = 0.01741933661286065
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 42.54352433333333
$ws.Range("N4").Value = 127.630573
$ws.Range("O4").Value = 0.1991972037564051
$ws.Range("P4").Value = 0.1991972037564051
$ws.Range("Q4").Value = 38.56230132622
$ws.Range("R4").Value = 347.06071193598
$ws.Range("S4").Value = 0.003469883144573411
$ws.Range("T4").Value = 0.00346988314457341

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Rarres2"
$ws.Range("C5").Value = "Ccrl2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.9064200000000001
$ws.Range("H5").Value = 2.71926
$ws.Range("I5").Value = 0.01741933661286065
$ws.Range("J5").Value = 0.01741933661286065
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.9041133333333334
$ws.Range("N5").Value = 2.71234
$ws.Range("O5").Value = 0.004233237624316298
$ws.Range("P5").Value = 0.004233237624316297
$ws.Range("Q5").Value = 0.8195064076000002
$ws.Range("R5").Value = 7.375557668400001
$ws.Range("S5").Value = [double]"7.374019114019214e-05"
$ws.Range("T5").Value = [double]"7.374019114019213e-05"

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Rarres2"
$ws.Range("C6").Value = "Ccrl2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 42.58841866666666
$ws.Range("H6").Value = 127.765256
$ws.Range("I6").Value = 0.8184528149909583
$ws.Range("J6").Value = 0.8184528149909585
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 169.8267096666667
$ws.Range("N6").Value = 509.480129
$ws.Range("O6").Value = 0.7951622771939805
$ws.Range("P6").Value = 0.7951622771939802
$ws.Range("Q6").Value = 7232.651012066446
$ws.Range("R6").Value = 65093.85910859802
$ws.Range("S6").Value = 0.650802804144034
$ws.Range("T6").Value = 0.6508028041440339

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Rarres2"
$ws.Range("C7").Value = "Ccrl2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 42.58841866666666
$ws.Range("H7").Value = 127.765256
$ws.Range("I7").Value = 0.8184528149909583
$ws.Range("J7").Value = 0.8184528149909585
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.30056
$ws.Range("N7").Value = 0.90168
$ws.Range("O7").Value = 0.001407281425298274
$ws.Range("P7").Value = 0.001407281425298273
$ws.Range("Q7").Value = 12.80037511445333
$ws.Range("R7").Value = 115.20337603008
$ws.Range("S7").Value = 0.00115179344401986
$ws.Range("T7").Value = 0.00115179344401986

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Rarres2"
$ws.Range("C8").Value = "Ccrl2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 42.58841866666666
$ws.Range("H8").Value = 127.765256
$ws.Range("I8").Value = 0.8184528149909583
$ws.Range("J8").Value = 0.8184528149909585
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 42.54352433333333
$ws.Range("N8").Value = 127.630573
$ws.Range("O8").Value = 0.1991972037564051
$ws.Range("P8").Value = 0.1991972037564051
$ws.Range("Q8").Value = 1811.86142586352
$ws.Range("R8").Value = 16306.75283277169
$ws.Range("S8").Value = 0.1630335121527572
$ws.Range("T8").Value = 0.1630335121527572

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Rarres2"
$ws.Range("C9").Value = "Ccrl2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 42.58841866666666
$ws.Range("H9").Value = 127.765256
$ws.Range("I9").Value = 0.8184528149909583
$ws.Range("J9").Value = 0.8184528149909585
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.9041133333333334
$ws.Range("N9").Value = 2.71234
$ws.Range("O9").Value = 0.004233237624316298
$ws.Range("P9").Value = 0.004233237624316297
$ws.Range("Q9").Value = 38.50475716211555
$ws.Range("R9").Value = 346.54281445904
$ws.Range("S9").Value = 0.003464705250147311
$ws.Range("T9").Value = 0.003464705250147311

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Rarres2"
$ws.Range("C10").Value = "Ccrl2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 8.540438
$ws.Range("H10").Value = 25.621314
$ws.Range("I10").Value = 0.164127848396181
$ws.Range("J10").Value = 0.164127848396181
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 169.8267096666667
$ws.Range("N10").Value = 509.480129
$ws.Range("O10").Value = 0.7951622771939805
$ws.Range("P10").Value = 0.7951622771939802
$ws.Range("Q10").Value = 1450.394484652167
$ws.Range("R10").Value = 13053.5503618695
$ws.Range("S10").Value = 0.1305082736816557
$ws.Range("T10").Value = 0.1305082736816556

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Rarres2"
$ws.Range("C11").Value = "Ccrl2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 8.540438
$ws.Range("H11").Value = 25.621314
$ws.Range("I11").Value = 0.164127848396181
$ws.Range("J11").Value = 0.164127848396181
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.30056
$ws.Range("N11").Value = 0.90168
$ws.Range("O11").Value = 0.001407281425298274
$ws.Range("P11").Value = 0.001407281425298273
$ws.Range("Q11").Value = 2.56691404528
$ws.Range("R11").Value = 23.10222640752
$ws.Range("S11").Value = 0.0002309740724221166
$ws.Range("T11").Value = 0.0002309740724221165

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Rarres2"
$ws.Range("C12").Value = "Ccrl2"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 8.540438
$ws.Range("H12").Value = 25.621314
$ws.Range("I12").Value = 0.164127848396181
$ws.Range("J12").Value = 0.164127848396181
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 42.54352433333333
$ws.Range("N12").Value = 127.630573
$ws.Range("O12").Value = 0.1991972037564051
$ws.Range("P12").Value = 0.1991972037564051
$ws.Range("Q12").Value = 363.3403318703246
$ws.Range("R12").Value = 3270.062986832922
$ws.Range("S12").Value = 0.03269380845907443
$ws.Range("T12").Value = 0.03269380845907443

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Rarres2"
$ws.Range("C13").Value = "Ccrl2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 8.540438
$ws.Range("H13").Value = 25.621314
$ws.Range("I13").Value = 0.164127848396181
$ws.Range("J13").Value = 0.164127848396181
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.9041133333333334
$ws.Range("N13").Value = 2.71234
$ws.Range("O13").Value = 0.004233237624316298
$ws.Range("P13").Value = 0.004233237624316297
$ws.Range("Q13").Value = 7.721523868306668
$ws.Range("R13").Value = 69.49371481476
$ws.Range("S13").Value = 0.0006947921830287948
$ws.Range("T13").Value = 0.0006947921830287947

Write-Host "Rebuilt LR-pair table. UsedRange: $($ws.UsedRange.Address())"